# Apply updated dSF (column F) values as per repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new dSF (column F) value
$updates = @{
    2  = -8
    5  = -4
    6  = -1
    7  = -6
    8  = -9
    10 = 8
    12 = -5
    14 = -6
    16 = -10
    17 = -6
    20 = -1
    21 = -3
    24 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
